# The dataset keeps its most-recent week at the top of the price history
# (rows 2..255). This weekly update inserts a brand-new record for the
# latest week right above the former row 239, pushing the existing rows
# 239..255 down to 240..256 (the sheet's dimension grows from T255 to T256).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("239:239").Insert()

$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44516
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = "Fruta"
$ws.Cells.Item(239, 7).Value = 100108
$ws.Cells.Item(239, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(239, 9).Value = 100108002
$ws.Cells.Item(239, 10).Value = "Mango"
$ws.Cells.Item(239, 11).Value = "Sin especificar"
$ws.Cells.Item(239, 12).Value = "Primera"
$ws.Cells.Item(239, 13).Value = 456
$ws.Cells.Item(239, 14).Value = 6000
$ws.Cells.Item(239, 15).Value = 6000
$ws.Cells.Item(239, 16).Value = 6000
$ws.Cells.Item(239, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(239, 18).Value = "Perú"
$ws.Cells.Item(239, 19).Value = 1500
$ws.Cells.Item(239, 20).Value = 4
